# save data done + era data updated
# Adds a new "Save" column (H) to the worksheet:
#  - H1 header "Save" with the same (bold/bordered) style as the other headers
#  - H2:H22 filled with 0/1 save indicator values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from G1 (bold font, border, centered/top alignment)
# onto H1, then set its text, matching the formatting used by the rest of
# the header row.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column values for each data row.
$saveValues = @(1, 1, 0, 0, 0, 0, 1, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 1, 0, 0)

$row = 2
foreach ($v in $saveValues) {
    $ws.Cells.Item($row, 8).Value = $v
    $row = $row + 1
}
